$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generated Date timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-10-03T16:37:46+01:00"

# --- Properties sheet: fill in the previously-empty "Uri" column ---
$props = $wb.Worksheets.Item("Properties")

# Row 2 -> "status" property
$props.Range("B2").Value = "http://hl7.org/fhir/concept-properties#status"

# Row 3 -> "severity" property
$props.Range("B3").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/CodeSystem/mobility-alert-level-cs#severity"
